$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text value into a cell without letting Excel's
# autoformatting reinterpret strings like "38%" as numbers (and without
# leaving a stray number-format behind on the cell).
function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2
Set-TextValue "C2" "18/05/2021 12:27:02"
Set-TextValue "D2" "38%"

# Row 3
Set-TextValue "C3" "19/05/2021 08:54:49"

# Row 4
Set-TextValue "C4" "19/05/2021 09:07:21"

# Row 6
Set-TextValue "C6" "19/05/2021 08:46:14"

# Row 7 - clear both cells
Set-TextValue "C7" ""
Set-TextValue "D7" ""

# Row 8
Set-TextValue "C8" "19/05/2021 09:01:22"
Set-TextValue "D8" "98%"

# Row 9
Set-TextValue "C9" "19/05/2021 09:14:27"
Set-TextValue "D9" "42%"

# Row 10
Set-TextValue "C10" "19/05/2021 08:57:16"

# Row 12
Set-TextValue "C12" "19/05/2021 09:12:19"

# Row 14
Set-TextValue "C14" "19/05/2021 08:49:14"

# Row 15
Set-TextValue "C15" "19/05/2021 08:49:43"

# Row 16
Set-TextValue "C16" "19/05/2021 09:16:58"

# Row 17
Set-TextValue "C17" "19/05/2021 09:14:33"
Set-TextValue "D17" "100%"

# Row 18 - clear both cells
Set-TextValue "C18" ""
Set-TextValue "D18" ""

# Row 20
Set-TextValue "C20" "19/05/2021 09:09:53"
Set-TextValue "D20" "90%"

# Row 21
Set-TextValue "C21" "19/05/2021 08:19:55"

# Row 22
Set-TextValue "C22" "18/05/2021 19:12:53"

# Row 23
Set-TextValue "C23" "19/05/2021 08:33:46"
Set-TextValue "D23" "72%"

# Row 24
Set-TextValue "C24" "19/05/2021 09:04:20"
Set-TextValue "D24" "86%"

# Row 25
Set-TextValue "C25" "19/05/2021 01:09:17"
Set-TextValue "D25" "28%"

# Row 27
Set-TextValue "C27" "19/05/2021 08:21:50"
Set-TextValue "D27" "36%"
